$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Paragraph 1: "Methods" -> "Data analysis". The "_GoBack" bookmark that
#    used to sit at the very start of the document is removed here (it gets
#    re-created later, right after the "- for Question 2: " run).
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(1)
$p1.Range.Find.Execute("Methods", $true, $false, $false, $false, $false, `
                        $true, 1, $false, "Data analysis", 2)

$d.Bookmarks.Item("_GoBack").Delete()

# ---------------------------------------------------------------------------
# 2) Paragraph 2 is currently empty; it gains the "What happened to
#    question 1?" note that used to live further down the document. A
#    single-character donor range (borrowed from the freshly renamed
#    paragraph 1, which carries exactly the run formatting every run in
#    this block shares) is temporarily overwritten with the desired text so
#    its FormattedText can be copied into place, then paragraph 1's text is
#    restored.
# ---------------------------------------------------------------------------
$p1b = $d.Paragraphs.Item(1)
$donor = $d.Range($p1b.Range.Start, $p1b.Range.Start + 1)
$donorFormatted = $donor.FormattedText
$donorFormatted.Text = "- What happened to question 1? "

$p2 = $d.Paragraphs.Item(2)
$insPoint = $d.Range($p2.Range.Start, $p2.Range.Start)
$insPoint.FormattedText = $donorFormatted

$p1c = $d.Paragraphs.Item(1)
$restoreRange = $d.Range($p1c.Range.Start, $p1c.Range.End - 1)
$restoreRange.Text = "Data analysis"

# ---------------------------------------------------------------------------
# 3) Paragraph 3 ("- It's definitely important ... broader population?")
#    loses all of its runs (and the gramStart/gramEnd proofing marks that
#    wrapped "definitely important"), becoming an empty paragraph.
# ---------------------------------------------------------------------------
$p3 = $d.Paragraphs.Item(3)
$clearRange = $d.Range($p3.Range.Start, $p3.Range.End - 1)
$clearRange.Delete()

# ---------------------------------------------------------------------------
# 4) The next seven paragraphs are removed outright:
#      - "You discuss this a bit ..."
#      - "This ties in to your proposed methods ..."
#      - "This is where it is also helpful ..."
#      - (empty)
#      - "Data analysis"
#      - "What happened to question 1?"
#      - (empty)
#    leaving "- for Question 2: " immediately after the now-empty paragraph.
# ---------------------------------------------------------------------------
$firstGone = $d.Paragraphs.Item(4)
$lastGone = $d.Paragraphs.Item(10)
$killRange = $d.Range($firstGone.Range.Start, $lastGone.Range.End)
$killRange.Delete()

# ---------------------------------------------------------------------------
# 5) Re-insert the "_GoBack" bookmark right after "- for Question 2: ".
#    A collapsed bookmark placed exactly at a paragraph's last character
#    position mis-resolves in this host, so a temporary placeholder
#    character is typed, the bookmark is anchored just before it, and the
#    placeholder is removed again.
# ---------------------------------------------------------------------------
$q2 = $d.Paragraphs.Item(4)
$tail = $d.Range($q2.Range.End - 1, $q2.Range.End - 1)
$tail.InsertAfter("#")

$q2b = $d.Paragraphs.Item(4)
$bmPoint = $d.Range($q2b.Range.End - 2, $q2b.Range.End - 2)
$d.Bookmarks.Add("_GoBack", $bmPoint)

$q2c = $d.Paragraphs.Item(4)
$placeholder = $d.Range($q2c.Range.End - 2, $q2c.Range.End - 1)
$placeholder.Delete()
